$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New aspect row: "city" — added at the bottom of the lookup table (row 31)
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "city"
$ws.Range("C31").Value = "City where process or stock is located, flows start or end"
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = "y"
$ws.Range("F31").Value = "Cit_y"

# Match the bold / centered formatting used by the other "aspect" (B) and
# "index_letter" (E) cells in this table.
foreach ($addr in @("B31", "E31")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Color = 0
    $cell.HorizontalAlignment = -4108
}

# Leave the same cell selected/active as in the authored workbook.
$ws.Range("E33").Select() | Out-Null
